$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 454) holds date-serial values that were all
# bumped from 45178 (2023-09-09) to 45179 (2023-09-10).
$ws.Range("C2:C454").Value = 45179
